$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3195.5
$ws.Range("I74").Value = 2679.111
$ws.Range("J74").Value = 3859.4285
$ws.Range("K74").Value = 2679.111
$ws.Range("L74").Value = 3859.4285
$ws.Range("M74").Value = -1743.111
$ws.Range("N74").Value = -5731.4285
$ws.Range("H76").Value = 23819176
$ws.Range("I76").Value = 13658.556
$ws.Range("J76").Value = 66669108
$ws.Range("K76").Value = 13658.556
$ws.Range("L76").Value = 66669108
$ws.Range("M76").Value = -13343.556
$ws.Range("N76").Value = -66669738
$ws.Range("H77").Value = 3195.5
$ws.Range("I77").Value = 2679.111
$ws.Range("J77").Value = 3859.4285
$ws.Range("K77").Value = 13395.555
$ws.Range("L77").Value = 19297.1425
$ws.Range("M77").Value = -8715.555
$ws.Range("N77").Value = -28657.1425
$ws.Range("H79").Value = 23819176
$ws.Range("I79").Value = 13658.556
$ws.Range("J79").Value = 66669108
$ws.Range("K79").Value = 13658.556
$ws.Range("L79").Value = 66669108
$ws.Range("M79").Value = -12566.556
$ws.Range("N79").Value = -66671292
$ws.Range("H112").Value = 1017.4915
$ws.Range("I112").Value = 525
$ws.Range("J112").Value = 1053.3091
$ws.Range("K112").Value = 1575
$ws.Range("L112").Value = 3159.9273
$ws.Range("M112").Value = -467
$ws.Range("N112").Value = -5375.927299999999
$ws.Range("H132").Value = 5003116.5
$ws.Range("I132").Value = 7045142
$ws.Range("J132").Value = 3674.2068
$ws.Range("K132").Value = 21135426
$ws.Range("L132").Value = 11022.6204
$ws.Range("M132").Value = -21132896
$ws.Range("N132").Value = -16082.6204
$ws.Range("H137").Value = 1318.6774
$ws.Range("I137").Value = 1194.9584
$ws.Range("J137").Value = 1742.8572
$ws.Range("K137").Value = 3584.8752
$ws.Range("L137").Value = 5228.571599999999
$ws.Range("M137").Value = -1034.8752
$ws.Range("N137").Value = -10328.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1255.9556
$ws.Range("I61").Value = 906.05554
$ws.Range("J61").Value = 2655.5557
$ws.Range("K61").Value = 906.05554
$ws.Range("L61").Value = 2655.5557
$ws.Range("M61").Value = -694.05554
$ws.Range("N61").Value = -3079.5557
$ws.Range("H62").Value = 19750
$ws.Range("J62").Value = 19750
$ws.Range("L62").Value = 19750
$ws.Range("N62").Value = -20998
$ws.Range("H65").Value = 19750
$ws.Range("J65").Value = 19750
$ws.Range("L65").Value = 59250
$ws.Range("N65").Value = -65490
$ws.Range("H74").Value = 1066.5264
$ws.Range("I74").Value = 1131.4193
$ws.Range("J74").Value = 779.1429000000001
$ws.Range("K74").Value = 1131.4193
$ws.Range("L74").Value = 779.1429000000001
$ws.Range("M74").Value = -257.4193
$ws.Range("N74").Value = -2527.1429
$ws.Range("H77").Value = 1066.5264
$ws.Range("I77").Value = 1131.4193
$ws.Range("J77").Value = 779.1429000000001
$ws.Range("K77").Value = 5657.0965
$ws.Range("L77").Value = 3895.7145
$ws.Range("M77").Value = -1289.0965
$ws.Range("N77").Value = -12631.7145
$ws.Range("H94").Value = 9310
$ws.Range("J94").Value = 9310
$ws.Range("L94").Value = 9310
$ws.Range("N94").Value = -11112
$ws.Range("H136").Value = 1255.9556
$ws.Range("I136").Value = 906.05554
$ws.Range("J136").Value = 2655.5557
$ws.Range("K136").Value = 2718.16662
$ws.Range("L136").Value = 7966.6671
$ws.Range("M136").Value = -168.16662
$ws.Range("N136").Value = -13066.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2225139.2
$ws.Range("I134").Value = 862.4358999999999
$ws.Range("J134").Value = 10111211
$ws.Range("K134").Value = 2587.3077
$ws.Range("L134").Value = 30333633
$ws.Range("M134").Value = -52.30769999999984
$ws.Range("N134").Value = -30338703

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 12000
$ws.Range("J4").Value = 12000
$ws.Range("L4").Value = 12000
$ws.Range("N4").Value = -12224
$ws.Range("H31").Value = 1285.5098
$ws.Range("I31").Value = 961.3043
$ws.Range("J31").Value = 1551.8214
$ws.Range("K31").Value = 961.3043
$ws.Range("L31").Value = 1551.8214
$ws.Range("M31").Value = -666.3043
$ws.Range("N31").Value = -2141.8214
$ws.Range("H34").Value = 1285.5098
$ws.Range("I34").Value = 961.3043
$ws.Range("J34").Value = 1551.8214
$ws.Range("K34").Value = 961.3043
$ws.Range("L34").Value = 1551.8214
$ws.Range("M34").Value = -759.3043
$ws.Range("N34").Value = -1955.8214
$ws.Range("H58").Value = 16950234
$ws.Range("I58").Value = 22728192
$ws.Range("J58").Value = 1552.9333
$ws.Range("K58").Value = 22728192
$ws.Range("L58").Value = 1552.9333
$ws.Range("M58").Value = -22727989
$ws.Range("N58").Value = -1958.9333
$ws.Range("H108").Value = 60000
$ws.Range("J108").Value = 60000
$ws.Range("L108").Value = 60000
$ws.Range("N108").Value = -67680
$ws.Range("H132").Value = 6945832
$ws.Range("I132").Value = 1056.4517
$ws.Range("J132").Value = 19609834
$ws.Range("K132").Value = 3169.3551
$ws.Range("L132").Value = 58829502
$ws.Range("M132").Value = -639.3551000000002
$ws.Range("N132").Value = -58834562
$ws.Range("H134").Value = 13889688
$ws.Range("I134").Value = 771.9375
$ws.Range("J134").Value = 125001016
$ws.Range("K134").Value = 2315.8125
$ws.Range("L134").Value = 375003048
$ws.Range("M134").Value = 219.1875
$ws.Range("N134").Value = -375008118
$ws.Range("H136").Value = 16950234
$ws.Range("I136").Value = 22728192
$ws.Range("J136").Value = 1552.9333
$ws.Range("K136").Value = 68184576
$ws.Range("L136").Value = 4658.7999
$ws.Range("M136").Value = -68182026
$ws.Range("N136").Value = -9758.7999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 166905
$ws.Range("I4").Value = 200268
$ws.Range("K4").Value = 600804
$ws.Range("M4").Value = -600692

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4817.7144
$ws.Range("I70").Value = 4706.357
$ws.Range("J70").Value = 4929.0713
$ws.Range("K70").Value = 4706.357
$ws.Range("L70").Value = 4929.0713
$ws.Range("M70").Value = -4436.357
$ws.Range("N70").Value = -5469.0713
$ws.Range("H73").Value = 4817.7144
$ws.Range("I73").Value = 4706.357
$ws.Range("J73").Value = 4929.0713
$ws.Range("K73").Value = 4706.357
$ws.Range("L73").Value = 4929.0713
$ws.Range("M73").Value = -3770.357
$ws.Range("N73").Value = -6801.0713
$ws.Range("H80").Value = 7696385
$ws.Range("I80").Value = 4800
$ws.Range("J80").Value = 33335002
$ws.Range("K80").Value = 4800
$ws.Range("L80").Value = 33335002
$ws.Range("M80").Value = -3802
$ws.Range("N80").Value = -33336998
$ws.Range("H83").Value = 7696385
$ws.Range("I83").Value = 4800
$ws.Range("J83").Value = 33335002
$ws.Range("K83").Value = 24000
$ws.Range("L83").Value = 166675010
$ws.Range("M83").Value = -19008
$ws.Range("N83").Value = -166684994
$ws.Range("H122").Value = 27784788
$ws.Range("I122").Value = 38470400
$ws.Range("K122").Value = 115411200
$ws.Range("M122").Value = -115408750
$ws.Range("H132").Value = 9601.111000000001
$ws.Range("I132").Value = 6678.1177
$ws.Range("J132").Value = 14570.2
$ws.Range("K132").Value = 20034.3531
$ws.Range("L132").Value = 43710.60000000001
$ws.Range("M132").Value = -17504.3531
$ws.Range("N132").Value = -48770.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 470000
$ws.Range("I2").Value = 1625000
$ws.Range("J2").Value = 50000
$ws.Range("K2").Value = 1625000
$ws.Range("L2").Value = 50000
$ws.Range("M2").Value = -1624888
$ws.Range("N2").Value = -50224
$ws.Range("H82").Value = 1751.5
$ws.Range("I82").Value = 1500
$ws.Range("J82").Value = 1835.3334
$ws.Range("K82").Value = 1500
$ws.Range("L82").Value = 1835.3334
$ws.Range("M82").Value = -1139
$ws.Range("N82").Value = -2557.3334
$ws.Range("H85").Value = 1751.5
$ws.Range("I85").Value = 1500
$ws.Range("J85").Value = 1835.3334
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 1835.3334
$ws.Range("M85").Value = -252
$ws.Range("N85").Value = -4331.3334
$ws.Range("H122").Value = 12640
$ws.Range("I122").Value = 14687.5
$ws.Range("J122").Value = 4450
$ws.Range("K122").Value = 44062.5
$ws.Range("L122").Value = 13350
$ws.Range("M122").Value = -41612.5
$ws.Range("N122").Value = -18250
$ws.Range("H132").Value = 27035258
$ws.Range("I132").Value = 52634188
$ws.Range("J132").Value = 14164.611
$ws.Range("K132").Value = 157902564
$ws.Range("L132").Value = 42493.833
$ws.Range("M132").Value = -157900034
$ws.Range("N132").Value = -47553.833
$ws.Range("H136").Value = 28946334
$ws.Range("I136").Value = 9922016
$ws.Range("J136").Value = 333335420
$ws.Range("K136").Value = 29766048
$ws.Range("L136").Value = 1000006260
$ws.Range("M136").Value = -29763498
$ws.Range("N136").Value = -1000011360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 20000
$ws.Range("I2").Value = 50000
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 50000
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -49888
$ws.Range("N2").Value = -5224
$ws.Range("H132").Value = 16723.717
$ws.Range("I132").Value = 21781.44
$ws.Range("J132").Value = 6186.7915
$ws.Range("K132").Value = 65344.31999999999
$ws.Range("L132").Value = 18560.3745
$ws.Range("M132").Value = -62814.31999999999
$ws.Range("N132").Value = -23620.3745
$ws.Range("H136").Value = 12821979
$ws.Range("I136").Value = 16129840
$ws.Range("J136").Value = 4016.25
$ws.Range("K136").Value = 48389520
$ws.Range("L136").Value = 12048.75
$ws.Range("M136").Value = -48386970
$ws.Range("N136").Value = -17148.75
